{"js": "// Append the new \"Results\" discussion paragraph text.\n// The target paragraph is the first (and only) empty paragraph that\n// immediately follows the bold \"Results\" heading; it currently contains\n// just a single tab character run.\nconst body = context.document.body;\n\n// Locate the \"Results\" heading (whole-word, case-sensitive match so we\n// don't match the lowercase \"results\" that appears earlier in the report).\nconst headingResults = body.search(\"Results\", { matchCase: true, matchWholeWord: true });\nheadingResults.load(\"items\");\nawait context.sync();\n\nconst heading = headingResults.items[0].paragraphs.getFirst();\nconst targetParagraph = heading.getNext();\n\nconst newText = \"The program works successfully, in so far as that the Search Bot responds to searches quickly, reliably, and accurately. The results are quick in that the Search Bot responds to user searches within one second. Faster responses would be unfeasible because there are limitations on how frequently the GroupMe API can be queried. Responding faster would require making API calls more frequently, which will result in the account being timed out. The results are reliable, in that the Search Bot responds to all messages in testing. To prevent the program from crashing, API calls are made within try-except statements, such that if the GroupMe API fails to return a response, the program continues running. If the program terminates, it will be able to respond to messages it missed upon being rebooted. Lastly, the program is run on a virtual machine in Azure rather than on one of our own computers, such that it is very unlikely its server will go down given the high reliability of Azure. Lastly, the results are accurate in that all of the messages returned match at least one keyword, and no documents are missed. However, there is a subjective tradeoff between keyword relevance and recency. If the search returns messages with the highest keyword relevance, through a ranking such as tf-idf, the results may likely be too old to still be relevant to the user. Thus, documents are returned in order of recency rather than pure keyword relevance. This could be an issue if a high percentage of documents were matches, but based on user testing, this is not the case. Since only a small percentage of documents are matches, it makes more sense to return matches in order of recency instead of relevance, in order to provide useful results. Based on user testing in a few different group chats, the results returned by the Search Bot are mostly relevant, and all relevant results are returned if they are among the most recent results. The GroupMe Search Bot is expected to provide a robust implementation of an important functionality missing in a popular messaging app.\";\n\n// Insert at the end of the paragraph (after the existing tab) so the\n// paragraph becomes: \"<tab>\" + newText, immediately before the\n// _GoBack bookmark that already closes out the paragraph.\ntargetParagraph.insertText(newText, Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Append the new \"Results\" discussion paragraph text.\n# The target paragraph is the first (and only) empty paragraph that\n# immediately follows the bold \"Results\" heading; it currently contains\n# just a single tab character run.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Results\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n$found = $rng.Find.Execute()\n\n$heading = $rng.Paragraphs(1)\n$targetParagraph = $heading.Next()\n\n$newText = 'The program works successfully, in so far as that the Search Bot responds to searches quickly, reliably, and accurately. The results are quick in that the Search Bot responds to user searches within one second. Faster responses would be unfeasible because there are limitations on how frequently the GroupMe API can be queried. Responding faster would require making API calls more frequently, which will result in the account being timed out. The results are reliable, in that the Search Bot responds to all messages in testing. To prevent the program from crashing, API calls are made within try-except statements, such that if the GroupMe API fails to return a response, the program continues running. If the program terminates, it will be able to respond to messages it missed upon being rebooted. Lastly, the program is run on a virtual machine in Azure rather than on one of our own computers, such that it is very unlikely its server will go down given the high reliability of Azure. Lastly, the results are accurate in that all of the messages returned match at least one keyword, and no documents are missed. However, there is a subjective tradeoff between keyword relevance and recency. If the search returns messages with the highest keyword relevance, through a ranking such as tf-idf, the results may likely be too old to still be relevant to the user. Thus, documents are returned in order of recency rather than pure keyword relevance. This could be an issue if a high percentage of documents were matches, but based on user testing, this is not the case. Since only a small percentage of documents are matches, it makes more sense to return matches in order of recency instead of relevance, in order to provide useful results. Based on user testing in a few different group chats, the results returned by the Search Bot are mostly relevant, and all relevant results are returned if they are among the most recent results. The GroupMe Search Bot is expected to provide a robust implementation of an important functionality missing in a popular messaging app.'\n\n# Append at the end of the paragraph (after the existing tab), leaving the\n# trailing _GoBack bookmark that already closes out the paragraph intact.\n$targetParagraph.Range.InsertAfter($newText)\n"}
